# Inserts a new weekly price record for "Brócoli" (Feria Lagunitas de Puerto
# Montt) before the existing row 481. All rows from 481 onward shift down by
# one (dimension grows from A1:R592 to A1:R593); the new row carries the same
# fixed attributes as the (now shifted) following row, with its own date and
# volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift row 481..592 down to 482..593, leaving a blank row 481 to fill in.
$ws.Rows.Item(481).Insert()

$newRow = 481
$ws.Cells.Item($newRow, 1).Value  = 4
$ws.Cells.Item($newRow, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($newRow, 3).Value  = "Los Lagos"
$ws.Cells.Item($newRow, 4).Value  = 45173
$ws.Cells.Item($newRow, 5).Value  = 10
$ws.Cells.Item($newRow, 6).Value  = 100112023
$ws.Cells.Item($newRow, 7).Value  = "Brócoli"
$ws.Cells.Item($newRow, 8).Value  = "Sin especificar"
$ws.Cells.Item($newRow, 9).Value  = "Primera"
$ws.Cells.Item($newRow, 10).Value = 500
$ws.Cells.Item($newRow, 11).Value = 1500
$ws.Cells.Item($newRow, 12).Value = 1500
$ws.Cells.Item($newRow, 13).Value = 1500
$ws.Cells.Item($newRow, 14).Value = "$/unidad"
$ws.Cells.Item($newRow, 15).Value = "Región Metropolitana"
$ws.Cells.Item($newRow, 16).Value = 1500
$ws.Cells.Item($newRow, 17).Value = 1
$ws.Cells.Item($newRow, 18).Value = "Hortaliza"
